# Mise a jour de l'application
# Fill in the "MG (%)" (body-fat %) column D for the newly-entered
# weigh-in batch (rows 133-159, date 45951) and move the view/selection
# the way the author left it after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# row -> MG (%) value (as a fraction, cell is formatted as 0.0%)
# rows not present here (137,139,143,151) only get the % number format
# applied, with no value - matching the source data set.
$values = [ordered]@{
    133 = 0.081
    134 = 0.108
    135 = 0.055
    136 = 0.037
    138 = 0.096
    140 = 0.037
    141 = 0.062
    142 = 0.093
    144 = 0.058
    145 = 0.062
    146 = 0.069
    147 = 0.058
    148 = 0.069
    149 = 0.051
    150 = 0.075
    152 = 0.083
    153 = 0.078
    154 = 0.086
    155 = 0.075
    156 = 0.083
    157 = 0.081
    158 = 0.062
    159 = 0.045
}

foreach ($row in 133..159) {
    $cell = $ws.Cells.Item($row, 4)
    if ($values.Contains($row)) {
        $cell.Value = $values[$row]
    }
    $cell.NumberFormat = "0.0%"
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4108
}

# Move the viewport / selection to match where the author ended up.
$ws.Range("A131").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 131
$win.ScrollColumn = 1
$ws.Range("G137").Select() | Out-Null
